$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) values are text that can look numeric (e.g. "242.11")
# or date-like (e.g. "29.842.10"); force text format first so Excel does not
# reinterpret them as numbers/dates, matching the original inline-string cells.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "29.842.10"
$ws.Cells.Item(3, 4).Value = "1.890.11"
$ws.Cells.Item(5, 4).Value = "0.7892"
$ws.Cells.Item(6, 4).Value = "242.11"
$ws.Cells.Item(8, 4).Value = "0.3187"
$ws.Cells.Item(9, 4).Value = "25.82"
$ws.Cells.Item(10, 4).Value = "0.07067"
$ws.Cells.Item(11, 4).Value = "0.08055"
$ws.Cells.Item(12, 4).Value = "0.7713"
$ws.Cells.Item(13, 4).Value = "1.907.99"
$ws.Cells.Item(14, 4).Value = "5.306"
$ws.Cells.Item(15, 4).Value = "92.26"
$ws.Cells.Item(16, 4).Value = "29.845.01"
$ws.Cells.Item(17, 4).Value = "13.86"
$ws.Cells.Item(18, 4).Value = "5.908"
$ws.Cells.Item(19, 4).Value = "243.65"
$ws.Cells.Item(20, 4).Value = "0.000007710"
$ws.Cells.Item(21, 4).Value = "1.001"
$ws.Cells.Item(22, 4).Value = "2.147.00"
$ws.Cells.Item(23, 4).Value = "8.048"
$ws.Cells.Item(25, 4).Value = "0.1621"
$ws.Cells.Item(26, 4).Value = "9.304"
$ws.Cells.Item(27, 4).Value = "164.89"
$ws.Cells.Item(29, 4).Value = "2.061"
$ws.Cells.Item(30, 4).Value = "1.379"
$ws.Cells.Item(31, 4).Value = "1.534"
$ws.Cells.Item(32, 4).Value = "4.402"
$ws.Cells.Item(33, 4).Value = "0.05626"
$ws.Cells.Item(34, 4).Value = "4.102"
$ws.Cells.Item(36, 4).Value = "0.7363"
$ws.Cells.Item(37, 4).Value = "1.002"
$ws.Cells.Item(38, 4).Value = "2.705"
$ws.Cells.Item(39, 4).Value = "0.01929"
$ws.Cells.Item(40, 4).Value = "2.772"
$ws.Cells.Item(41, 4).Value = "0.4446"
$ws.Cells.Item(42, 4).Value = "72.22"
$ws.Cells.Item(43, 4).Value = "5.865"
$ws.Cells.Item(44, 4).Value = "0.8455"
$ws.Cells.Item(45, 4).Value = "1.001"
$ws.Cells.Item(46, 4).Value = "1.882"
$ws.Cells.Item(47, 4).Value = "102.41"
$ws.Cells.Item(48, 4).Value = "1.020.52"
$ws.Cells.Item(49, 4).Value = "9.948"
$ws.Cells.Item(50, 4).Value = "7.478"
$ws.Cells.Item(51, 4).Value = "2.958"

# Column E (Volume/1h) values always carry surrounding spaces and a "%" sign
# so Excel keeps them as plain text without any extra formatting needed.
$ws.Cells.Item(2, 5).Value = "  -0.27%  "
$ws.Cells.Item(3, 5).Value = "  -0.64%  "
$ws.Cells.Item(4, 5).Value = "  +0.23%  "
$ws.Cells.Item(5, 5).Value = "  -1.40%  "
$ws.Cells.Item(6, 5).Value = "  +0.52%  "
$ws.Cells.Item(7, 5).Value = "  +0.17%  "
$ws.Cells.Item(8, 5).Value = "  +2.05%  "
$ws.Cells.Item(9, 5).Value = "  -1.72%  "
$ws.Cells.Item(10, 5).Value = "  -0.06%  "
$ws.Cells.Item(11, 5).Value = "  +0.99%  "
$ws.Cells.Item(12, 5).Value = "  +4.57%  "
$ws.Cells.Item(13, 5).Value = "  +0.38%  "
$ws.Cells.Item(14, 5).Value = "  +2.46%  "
$ws.Cells.Item(15, 5).Value = "  -0.29%  "
$ws.Cells.Item(16, 5).Value = "  -0.23%  "
$ws.Cells.Item(17, 5).Value = "  -0.52%  "
$ws.Cells.Item(18, 5).Value = "  +0.67%  "
$ws.Cells.Item(19, 5).Value = "  -0.52%  "
$ws.Cells.Item(20, 5).Value = "  -0.67%  "
$ws.Cells.Item(21, 5).Value = "  +0.11%  "
$ws.Cells.Item(22, 5).Value = "  -0.14%  "
$ws.Cells.Item(23, 5).Value = "  +16.62%  "
$ws.Cells.Item(24, 5).Value = "  +0.20%  "
$ws.Cells.Item(25, 5).Value = "  +13.90%  "
$ws.Cells.Item(26, 5).Value = "  +1.23%  "
$ws.Cells.Item(27, 5).Value = "  -1.57%  "
$ws.Cells.Item(28, 5).Value = "  -0.72%  "
$ws.Cells.Item(29, 5).Value = "  +1.29%  "
$ws.Cells.Item(30, 5).Value = "  +1.73%  "
$ws.Cells.Item(31, 5).Value = "  +1.43%  "
$ws.Cells.Item(32, 5).Value = "  +2.45%  "
$ws.Cells.Item(33, 5).Value = "  +0.87%  "
$ws.Cells.Item(34, 5).Value = "  +1.03%  "
$ws.Cells.Item(35, 5).Value = "  +0.40%  "
$ws.Cells.Item(36, 5).Value = "  +1.00%  "
$ws.Cells.Item(37, 5).Value = "  +0.29%  "
$ws.Cells.Item(38, 5).Value = "  -0.46%  "
$ws.Cells.Item(39, 5).Value = "  -0.07%  "
$ws.Cells.Item(40, 5).Value = "  -0.43%  "
$ws.Cells.Item(41, 5).Value = "  +0.94%  "
$ws.Cells.Item(42, 5).Value = "  +0.25%  "
$ws.Cells.Item(43, 5).Value = "  -2.20%  "
$ws.Cells.Item(44, 5).Value = "  +1.20%  "
$ws.Cells.Item(45, 5).Value = "  +0.11%  "
$ws.Cells.Item(46, 5).Value = "  +0.91%  "
$ws.Cells.Item(47, 5).Value = "  +1.96%  "
$ws.Cells.Item(48, 5).Value = "  +4.39%  "
$ws.Cells.Item(49, 5).Value = "  +2.39%  "
$ws.Cells.Item(50, 5).Value = "  -1.18%  "
$ws.Cells.Item(51, 5).Value = "  +7.29%  "

# Restore the default "Normal" style on column D so no stray style index is
# introduced (the workbook keeps these cells on the default style).
$priceRange.Style = "Normal"
